$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.328.30"
$ws.Range("E2").Value = "  -4.64%  "
$ws.Range("D3").Value = "3.379.37"
$ws.Range("E3").Value = "  -6.73%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'186.72"
$ws.Range("E5").Value = "  -8.42%  "
$ws.Range("D6").Value = "'527.44"
$ws.Range("E6").Value = "  -7.26%  "
$ws.Range("D7").Value = "'0.606"
$ws.Range("E7").Value = "  -2.83%  "
$ws.Range("D8").Value = "3.377.38"
$ws.Range("E8").Value = "  -6.66%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'0.626"
$ws.Range("E10").Value = "  -7.69%  "
$ws.Range("D11").Value = "'59.01"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("D12").Value = "'0.133"
$ws.Range("E12").Value = "  -12.75%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -11.83%  "
$ws.Range("D14").Value = "'9.25"
$ws.Range("E14").Value = "  -8.52%  "
$ws.Range("D15").Value = "3.934.06"
$ws.Range("E15").Value = "  -6.30%  "
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("D17").Value = "3.383.89"
$ws.Range("E17").Value = "  -6.19%  "
$ws.Range("D18").Value = "65.154.31"
$ws.Range("E18").Value = "  -4.68%  "
$ws.Range("D19").Value = "'17.44"
$ws.Range("E19").Value = "  -8.67%  "
$ws.Range("D20").Value = "'11.13"
$ws.Range("E20").Value = "  -10.47%  "
$ws.Range("D21").Value = "'0.973"
$ws.Range("E21").Value = "  -9.97%  "
$ws.Range("D22").Value = "'372.61"
$ws.Range("E22").Value = "  -8.05%  "
$ws.Range("D23").Value = "'81.84"
$ws.Range("E23").Value = "  -4.50%  "
$ws.Range("D24").Value = "'3.72"
$ws.Range("E24").Value = "  -11.21%  "
$ws.Range("D25").Value = "'10.84"
$ws.Range("E25").Value = "  -17.30%  "
$ws.Range("D26").Value = "'3.70"
$ws.Range("E26").Value = "  -5.11%  "
$ws.Range("D27").Value = "'2.65"
$ws.Range("E27").Value = "  -9.95%  "
$ws.Range("D28").Value = "'11.47"
$ws.Range("E28").Value = "  -9.54%  "
$ws.Range("D29").Value = "'8.53"
$ws.Range("E29").Value = "  -9.52%  "
$ws.Range("D30").Value = "'686.78"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("D31").Value = "'29.63"
$ws.Range("E31").Value = "  -6.64%  "
$ws.Range("D32").Value = "'6.72"
$ws.Range("E32").Value = "  -18.07%  "
$ws.Range("D33").Value = "'11.17"
$ws.Range("E33").Value = "  -9.29%  "
$ws.Range("D34").Value = "'61.21"
$ws.Range("E34").Value = "  -4.32%  "
$ws.Range("D35").Value = "'0.106"
$ws.Range("E35").Value = "  -7.87%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "'36.52"
$ws.Range("E37").Value = "  -13.78%  "
$ws.Range("D38").Value = "'0.383"
$ws.Range("E38").Value = "  -9.60%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "'0.128"
$ws.Range("E40").Value = "  -6.46%  "
$ws.Range("D41").Value = "2.862.50"
$ws.Range("E41").Value = "  -12.73%  "
$ws.Range("D42").Value = "'2.77"
$ws.Range("E42").Value = "  -12.97%  "
$ws.Range("D43").Value = "'2.67"
$ws.Range("E43").Value = "  -4.56%  "
$ws.Range("D44").Value = "'0.0396"
$ws.Range("E44").Value = "  -5.81%  "
$ws.Range("D45").Value = "0.0₃0620"
$ws.Range("E45").Value = "  -20.50%  "
$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = "  -15.06%  "
$ws.Range("D47").Value = "'0.126"
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("D48").Value = "'137.13"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "'2.63"
$ws.Range("E49").Value = "  -4.73%  "
$ws.Range("D50").Value = "'2.84"
$ws.Range("E50").Value = "  -8.97%  "
$ws.Range("D51").Value = "'7.68"
$ws.Range("E51").Value = "  -14.01%  "
